$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Set D4 numeric mark value (this feeds formula in J4, recalculated automatically)
$ws.Range("D4").Value = 6

# Set D5 comment text
$ws.Range("D5").Value = "The D&C version is not implemented. Please, check the video of the last seminar (the implementation should be very similar to Mergesort) *** The problem with the time is that you are executing the algorithm in the Constructor but when you calculate the time (t2-t1) you are just calling start() that the only thing it does is to return the previously calculated value *** The rest of the things are fine"

# Update row 12 height
$ws.Rows.Item(12).RowHeight = 102.75

# Update selection to D5:D12
$ws.Range("D5:D12").Select()
